# Update the "score_at_condition" (column P) values on the
# Condition_by_Condition sheet, and the "overall_score" / "consistency"
# values (columns B, C) on the Overall_Rankings sheet, to reflect the
# new priority order used when computing these composite scores.

$wb = $excel.ActiveWorkbook

$wsCond = $wb.Worksheets.Item("Condition_by_Condition")
$wsRank = $wb.Worksheets.Item("Overall_Rankings")

# score_at_condition (column P), rows 2-46
$newScores = @{
    2  = 4.4
    3  = 3.37155336908946
    4  = 0.6241629712659377
    5  = 4.4
    6  = 3.347926146658295
    7  = 0.6241340219678418
    8  = 4.4
    9  = 3.336389906470643
    10 = 0.6240764017126385
    11 = 3.702032520325195
    12 = 2.758747493241533
    13 = 1.475247813654563
    14 = 3.730408163265326
    15 = 2.795013647031681
    16 = 1.475223093437327
    17 = 3.781504065040667
    18 = 2.867116127102249
    19 = 1.475186003396205
    20 = 3.9
    21 = 2.785151743897122
    22 = 1.475004882786253
    23 = 3.920833333333339
    24 = 2.824366798369596
    25 = 1.474991480101823
    26 = 3.951100244498751
    27 = 2.901803361749447
    28 = 1.474972978230662
    29 = 4.343333333333266
    30 = 2.683396196626719
    31 = 1.474617276804007
    32 = 4.329166666666734
    33 = 2.715338346136939
    34 = 1.474619950802749
    35 = 4.336250000000101
    36 = 2.778322100352302
    37 = 1.474618613795579
    38 = 4.4
    39 = 2.51506446418664
    40 = 1.474579848440692
    41 = 4.4
    42 = 2.536153784244761
    43 = 1.474579313735402
    44 = 4.4
    45 = 2.577869897105686
    46 = 1.47457958108827
}

foreach ($row in $newScores.Keys) {
    $wsCond.Range("P$row").Value = $newScores[$row]
}

# overall_score (column B) and consistency (column C), rows 2-4
$wsRank.Range("B2").Value = 4.159641888430892
$wsRank.Range("C2").Value = 0.7835335425850183

$wsRank.Range("B3").Value = 2.852947558817539
$wsRank.Range("C3").Value = 0.7860433831129233

$wsRank.Range("B4").Value = 1.30470628208133
$wsRank.Range("C4").Value = 0.746106633459375
